$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stale TestResult ("Fail") for the test case on row 3 (RunFlag = 0,
# so this case is not executed and should not carry a leftover result).
$ws.Range("C3").ClearContents()

# Move/restore the active selection, matching the saved workbook view state.
$ws.Range("C6").Select()
